$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 9 through 13 (entire rows), shifting data up not needed since rows 4-8 are updated in place
$ws.Range("A9:N13").EntireRow.Delete() | Out-Null

# Update rows 4-8 with new sample data
$ws.Range("A4").Value = "SS-1"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.45
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1.67
$ws.Range("F4").Value = 59.1
$ws.Range("G4").Value = 27.05
$ws.Range("N4").Value = 82

$ws.Range("A5").Value = "SS-2"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2.4500000000000002
$ws.Range("D5").Value = 1.02
$ws.Range("E5").Value = 1.83
$ws.Range("F5").Value = 77.01
$ws.Range("G5").Value = 27.55
$ws.Range("N5").Value = 60

$ws.Range("A6").Value = "SS-3"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3.45
$ws.Range("D6").Value = 3.73
$ws.Range("E6").Value = 1.9
$ws.Range("F6").Value = 93.99
$ws.Range("G6").Value = 33.87
$ws.Range("N6").Value = 75

$ws.Range("A7").Value = "SS-4"
$ws.Range("B7").Value = 4.5
$ws.Range("C7").Value = 4.95
$ws.Range("D7").Value = 3.84
$ws.Range("E7").Value = 1.9
$ws.Range("F7").Value = 48.09
$ws.Range("G7").Value = 15.35
$ws.Range("N7").Value = 80

$ws.Range("A8").Value = "SS-5"
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 6.45
$ws.Range("D8").Value = 3.34
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 63.22
$ws.Range("G8").Value = 17.3
$ws.Range("N8").Value = 100

$ws.Range("N9").Select() | Out-Null
